# feat: add 2022-Q4 data
#
# Summary of the change (see commit message "feat: add 2022-Q4 data"):
#   - "总计" row 2 ("2022-Q1"/0.02) is now row 3; row 2 instead reports
#     the brand-new "2022-Q4" figures (same count 2, value 0).
#   - The existing "2022-Q1" worksheet is repurposed in place (keeps its
#     sheetId/rId) to hold the "2022-Q4" fund table; a new worksheet
#     named "2022-Q1" is appended right after it with the original fund
#     rows restored.

# ------------------------------------------------------------------
# Helper: write a *text* value into a cell without Excel's automatic
# "looks like a number" coercion (e.g. "0.21" / "015921" must stay
# text, matching the source file, not become 0.21 / 15921). We do this
# by writing a formula that RETURNS the literal string into a scratch
# cell, then Copy + PasteSpecial(Values) into the destination -- that
# carries the text type over without minting any new cell style.
# ------------------------------------------------------------------
function Set-TextCell {
    param($ws, $addr, $text)
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace("""", """""")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
}

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)     # "总计"
$q1Sheet = $wb.Worksheets.Item(2)     # existing "2022-Q1" sheet (sheetId 2 / rId2)

# ------------------------------------------------------------------
# 1) Known original "2022-Q1" fund-holding data (read back from the
#    source workbook) -- needed once we overwrite $q1Sheet in place.
# ------------------------------------------------------------------
$oldCode1   = "620007"
$oldName1   = "金元顺安优质精选灵活配置混合A"
$oldScale1  = "0.75"
$oldPos1    = "39.31"
$oldRatio1  = "1.43"
$oldValue1  = "0.0107"
$oldRank1   = 9

$oldCode2   = "001375"
$oldName2   = "金元顺安优质精选灵活配置混合C"
$oldScale2  = "0.69"
$oldPos2    = "39.31"
$oldRatio2  = "1.43"
$oldValue2  = "0.0099"
$oldRank2   = 9

# ------------------------------------------------------------------
# 2) Insert a brand-new worksheet right after the existing "2022-Q1"
#    tab; it will hold the *restored* "2022-Q1" figures. Copy that
#    sheet's current header/A-column formatting onto it first, while
#    the source still carries it. The existing tab is renamed out of
#    the way first so the new tab can claim the "2022-Q1" name.
# ------------------------------------------------------------------
$q1Sheet.Name = "2022-Q1-old"
$newQ1 = $wb.Worksheets.Add($null, $q1Sheet)
$newQ1.Name = "2022-Q1"

$q1Sheet.Range("B1:H1").Copy() | Out-Null
$newQ1.Range("B1").PasteSpecial(-4122) | Out-Null
$q1Sheet.Range("A2:A3").Copy() | Out-Null
$newQ1.Range("A2").PasteSpecial(-4122) | Out-Null

Set-TextCell $newQ1 "B1" "基金代码"
Set-TextCell $newQ1 "C1" "基金名称"
Set-TextCell $newQ1 "D1" "基金规模"
Set-TextCell $newQ1 "E1" "股票总仓位"
Set-TextCell $newQ1 "F1" "仓位占比"
Set-TextCell $newQ1 "G1" "持有市值(亿元)"
Set-TextCell $newQ1 "H1" "仓位排名"

$newQ1.Range("A2").Value = 0
Set-TextCell $newQ1 "B2" $oldCode1
Set-TextCell $newQ1 "C2" $oldName1
Set-TextCell $newQ1 "D2" $oldScale1
Set-TextCell $newQ1 "E2" $oldPos1
Set-TextCell $newQ1 "F2" $oldRatio1
Set-TextCell $newQ1 "G2" $oldValue1
$newQ1.Range("H2").Value = $oldRank1

$newQ1.Range("A3").Value = 1
Set-TextCell $newQ1 "B3" $oldCode2
Set-TextCell $newQ1 "C3" $oldName2
Set-TextCell $newQ1 "D3" $oldScale2
Set-TextCell $newQ1 "E3" $oldPos2
Set-TextCell $newQ1 "F3" $oldRatio2
Set-TextCell $newQ1 "G3" $oldValue2
$newQ1.Range("H3").Value = $oldRank2

# ------------------------------------------------------------------
# 3) Repurpose the original worksheet (keeps its sheetId/rId) as the
#    new "2022-Q4" sheet: rename it, restyle its header/A-column to
#    match "总计" (reusing that style index instead of minting a new
#    one), and overwrite its fund-holding rows with the Q4 figures.
# ------------------------------------------------------------------
$q1Sheet.Name = "2022-Q4"

$summary.Range("B1:D1").Copy() | Out-Null
$q1Sheet.Range("B1").PasteSpecial(-4122) | Out-Null
$q1Sheet.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$summary.Range("A2").Copy() | Out-Null
$q1Sheet.Range("A2").PasteSpecial(-4122) | Out-Null
$q1Sheet.Range("A3").PasteSpecial(-4122) | Out-Null

Set-TextCell $q1Sheet "B1" "基金代码"
Set-TextCell $q1Sheet "C1" "基金名称"
Set-TextCell $q1Sheet "D1" "基金规模"
Set-TextCell $q1Sheet "E1" "股票总仓位"
Set-TextCell $q1Sheet "F1" "仓位占比"
Set-TextCell $q1Sheet "G1" "持有市值(亿元)"
Set-TextCell $q1Sheet "H1" "仓位排名"

$q1Sheet.Range("A2").Value = 0
Set-TextCell $q1Sheet "B2" "015921"
Set-TextCell $q1Sheet "C2" "申万菱信国证2000指数增强A"
Set-TextCell $q1Sheet "D2" "0.21"
Set-TextCell $q1Sheet "E2" "94.00"
Set-TextCell $q1Sheet "F2" "0.52"
Set-TextCell $q1Sheet "G2" "0.0011"
$q1Sheet.Range("H2").Value = 4

$q1Sheet.Range("A3").Value = 1
Set-TextCell $q1Sheet "B3" "015922"
Set-TextCell $q1Sheet "C3" "申万菱信国证2000指数增强C"
Set-TextCell $q1Sheet "D3" "0.08"
Set-TextCell $q1Sheet "E3" "94.00"
Set-TextCell $q1Sheet "F3" "0.52"
Set-TextCell $q1Sheet "G3" "0.0004"
$q1Sheet.Range("H3").Value = 4

# ------------------------------------------------------------------
# 4) "总计" sheet: row 2 used to describe "2022-Q1"; it now describes
#    "2022-Q4" (same count, zero value-change), and a new row 3 is
#    appended carrying what used to be row 2's "2022-Q1" data.
# ------------------------------------------------------------------
$summary.Range("A2:D2").Copy() | Out-Null
$summary.Range("A3").PasteSpecial(-4122) | Out-Null

$summary.Range("A3").Value = 1
Set-TextCell $summary "B3" "2022-Q1"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.02

Set-TextCell $summary "B2" "2022-Q4"
$summary.Range("D2").Value = 0

# ------------------------------------------------------------------
# 5) Leave the workbook selection/active sheet on "总计", matching the
#    unchanged <bookViews> (activeTab stays 0).
# ------------------------------------------------------------------
$summary.Activate()
$summary.Range("A1").Select() | Out-Null
